$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new precondition bullet after "Es müssen bereits
#    Nährstoffe im System vorhanden sein." :
#    "Dialog für das erstellen ein Rezept ist bereits geöffnet.(TeamC)"
# ------------------------------------------------------------------
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Es müssen bereits Nährstoffe im System vorhanden sein.", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$anchorPara = $anchor.Paragraphs(1)

# Adding a paragraph after it inherits the same list (numId/ilvl) and
# run formatting as the anchor paragraph.
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()
$newPara.Range.Text = "Dialog für das erstellen ein Rezept ist bereits geöffnet.(TeamC)"

# ------------------------------------------------------------------
# 2) Remove the "Zubereitung (Pflichtfeld)" bullet (together with its
#    paragraph mark / _GoBack bookmark) that used to follow
#    "Zubereitungszeit (Pflichtfeld)".
# ------------------------------------------------------------------
$target = $d.Content.Duplicate
$target.Find.Execute("Zubereitung (Pflichtfeld)", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$targetPara = $target.Paragraphs(1)
$targetPara.Range.Delete()
